$wb = $excel.ActiveWorkbook

# Overview sheet: rows for 289a79fa (row 3) and 37da4d39 (row 4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-18 14:16:24"
$wsOverview.Range("G4").Value = "2016-08-18 14:16:24"

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-18 14:16:19"
$wsZhCn.Range("H4").Value = "2016-08-18 14:16:19"
$wsZhCn.Range("K3").Value = "2016-08-18 14:16:35"
$wsZhCn.Range("K4").Value = "2016-08-18 14:16:35"

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-18 14:16:24"
$wsDeDe.Range("H4").Value = "2016-08-18 14:16:24"
$wsDeDe.Range("K3").Value = "2016-08-18 14:16:43"
$wsDeDe.Range("K4").Value = "2016-08-18 14:16:43"
